$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Insert a brand-new row above the existing row 70 ("Progression from CD4>500 to
# CD4>350 on unsuppressive ART" / usvlproggt500), pushing it and everything below
# down by one. This mirrors the author's insertion of a new "Number of VL tests
# recommended per person per year" / requiredvl parameter row.
$ws.Rows("70").Insert()

# Populate the newly-opened row 70 with the new parameter's data.
$ws.Range("A70").Value = "Number of VL tests recommended per person per year"
$ws.Range("B70").Value = "constant"
$ws.Range("C70").Value = "requiredvl"
$ws.Range("D70").Value = "requiredvl"
$ws.Range("E70").Value = "(0, 'maxacts')"
$ws.Range("F70").Value = "tot"
$ws.Range("G70").Value = "constant"
$ws.Range("H70").Value = "const"
$ws.Range("I70").Value = "None"
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = "None"
$ws.Range("L70").Value = 1
$ws.Range("M70").Value = "const"

# Match the author's resulting selection (F70:M70, active cell F70).
$ws.Range("F70:M70").Select()
